# Apply updated cryptocurrency price/volume data to the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.482.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.923.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.73%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.695"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.90"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "57.89"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.365"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0766"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0995"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.799"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.94%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.193.80"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.58%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.914.27"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.20%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.420.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0844"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "252.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.58%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.04%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.38%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.84%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.85"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.85"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.91%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.81"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.37%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.72%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.22%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.89%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.96%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0844"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +21.67%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -15.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.864"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.23%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.48%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "104.72"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0230"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.28%  "

# Row 42
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.25"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.40%  "

# Row 43
$ws.Range("B43").Value = "Gas"
$ws.Range("C43").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +22.83%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.347.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.80%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.40"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0810"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.83%  "

# Row 48
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.87%  "

# Row 49
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.44"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.96%  "
